$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2 updates
$ws.Range("H2").Value = 6.2
$ws.Range("J2").Value = 3.85
$ws.Range("N2").Value = 4.1
$ws.Range("P2").Value = 2.08
$ws.Range("S2").Value = 2.88
$ws.Range("X2").Value = 1000
$ws.Range("AG2").Value = 12

# Row 3 updates
$ws.Range("F3").Value = 8
$ws.Range("G3").Value = 9.800000000000001
$ws.Range("I3").Value = 1.45
$ws.Range("P3").Value = 2.34
$ws.Range("R3").Value = 1.53
$ws.Range("S3").Value = 2.5
$ws.Range("U3").Value = 1.94
$ws.Range("V3").Value = 3.2
$ws.Range("AG3").Value = 34
$ws.Range("AH3").Value = 1000

# Row 4 updates
$ws.Range("F4").Value = 1.09
$ws.Range("H4").Value = 1.09
$ws.Range("J4").Value = 1.17
$ws.Range("Q4").Value = 1.2
$ws.Range("S4").Value = 1.21
